$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# 1) e007 Morning Briefing - Weather Roll (row 8): drop trailing "B24"
$ws.Range("B8").Value = "<Bold>e007 Morning Briefing - Weather Roll</Bold> <InlineUIContainer><Button Content='r4.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`n<LineBreak/><LineBreak/>`nThe `n<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n Table determines weather for today:  `n<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>`n<LineBreak/>"

# 2) e009 Ammo Loading Limits (row 10): renamed from e008 -> e009, reformatted AP/HE lines
$ws.Range("B10").Value = "<Bold>e009 Ammo Loading Limits</Bold> <InlineUIContainer><Button Content='r16.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`n<LineBreak/><LineBreak/>`nSee `n<InlineUIContainer><Button Content='r16.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nfor ammo types. See `n<InlineUIContainer><Button Content='r16.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nfor loading ammo. The Tank Card limits the number of normal main gun ammo allowed to AMMO_NORMAL_LOAD. Extra ammo is added in a later step `n<InlineUIContainer><Button Content='e009b' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.`n<LineBreak/><LineBreak/>`n <Bold>--AP:</Bold> Unlimited<LineBreak/>`n <Bold>-- HE:</Bold> Unlimited"

# 3) e010 Time Check (row 11): reflowed line breaks
$ws.Range("B11").Value = "<Bold>e010 Time Check</Bold> `n<InlineUIContainer><Button Content='r4.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  `n<LineBreak/><LineBreak/>`nDetermine sunrise and sunset for current month using the <InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. `nRoll 1D/2  on the `n<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table. `nThe Time Table also provides the timed used for each action take. Additionally, the same die roll is used to determine the ammo expended:  `n<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>`n<LineBreak/><LineBreak/>"

# 4) Insert a new row 33 for e032 "No Combat" (pushes old rows 33-37 down to 34-38)
$ws.Rows.Item(33).Insert()
$ws.Range("A33").Value = "e032"
$ws.Range("B33").Value = "<Bold>e032 No Combat</Bold> `n<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nNo combat. Victory points added to the After Action Report `n<InlineUIContainer><Image Name='AAR' Height='60' Width='60'></Image></InlineUIContainer> `nto reflect area under US Control. Continue with "
$ws.Range("A33:B33").EntireRow.RowHeight = 90
$ws.Cells.Item(33, 1).Style = $ws.Cells.Item(32, 1).Style
$ws.Cells.Item(33, 2).Style = $ws.Cells.Item(32, 2).Style

# 5) e031 Resistance Table -> e031 Battle Check (row 32, still row 32 since insert was below it)
$ws.Range("B32").Value = "<Bold>e031 Battle Check</Bold> `n<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nRoll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nTable to determine if combat occurs in this area: <LineBreak/><LineBreak/>`nDie Roll =  <InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer> "

# Update selection to match final state
$ws.Range("B33").Select()
